# Update "ES" worksheet: several task statuses moved from "pm"/"ep" to "f"
# (finished), and the now-redundant "assigned to" cells next to them were
# cleared. Selection/active cell also moved from I15 to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ES")

# Status -> "f" (Feito/Done), and clear the adjacent "who" cell.
$ws.Range("D4").Value = "f"
$ws.Range("E4").Value = ""

$ws.Range("D6").Value = "f"
$ws.Range("E6").Value = ""

$ws.Range("B8").Value = "f"
$ws.Range("C8").Value = ""

$ws.Range("B9").Value = "f"
$ws.Range("C9").Value = ""

$ws.Range("B14").Value = "f"
$ws.Range("C14").Value = ""

$ws.Range("D14").Value = "f"
$ws.Range("E14").Value = ""

$ws.Range("F14").Value = "f"
$ws.Range("G14").Value = ""

$ws.Range("F15").Value = "f"
$ws.Range("G15").Value = ""

# Move the active selection to D7, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("D7").Select() | Out-Null
